$wb = $excel.ActiveWorkbook
$wsObras = $wb.Worksheets.Item("Obras")
$wsGas = $wb.Worksheets.Item("Gasto de gasolina")

# Copy the formatting (bold, wrap, vertical-center) from an already-styled
# header cell onto A2:A6 of "Gasto de gasolina" so it matches the style
# used for the same column in "Obras".
$wsGas.Range("B1").Copy()
$wsGas.Range("A2:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Correct / normalize the location names in column A so they match the
# same set of labels used in "Obras".
$wsGas.Range("A2").Value = "Oficina"
$wsGas.Range("A3").Value = "Naucalpan"
$wsGas.Range("A4").Value = "Santa fe"
$wsGas.Range("A5").Value = "Lomas"
$wsGas.Range("A6").Value = "Satelite"

# Move the active selection/sheet: "Obras" becomes the active tab with
# selection at I16; "Gasto de gasolina" selection moves to A6.
$wsGas.Range("A6").Select() | Out-Null
$wsObras.Activate() | Out-Null
$wsObras.Range("I16").Select() | Out-Null
